# Insert a new row at position 97 (shifts existing rows 97:134 down to 98:135)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(97).Insert()

# Populate the newly inserted row 97 with the new record
$ws.Cells.Item(97, 1).Value2 = 5
$ws.Cells.Item(97, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(97, 3).Value2 = "Maule"
$ws.Cells.Item(97, 4).Value2 = 45027
$ws.Cells.Item(97, 5).Value2 = 7
$ws.Cells.Item(97, 6).Value2 = "Fruta"
$ws.Cells.Item(97, 7).Value2 = 100101
$ws.Cells.Item(97, 8).Value2 = "Berries"
$ws.Cells.Item(97, 9).Value2 = 100101001
$ws.Cells.Item(97, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item(97, 11).Value2 = "Sin especificar"
$ws.Cells.Item(97, 12).Value2 = "Segunda"
$ws.Cells.Item(97, 13).Value2 = 50
$ws.Cells.Item(97, 14).Value2 = 4000
$ws.Cells.Item(97, 15).Value2 = 4000
$ws.Cells.Item(97, 16).Value2 = 4000
$ws.Cells.Item(97, 17).Value2 = "$/bandeja 2 kilos"
$ws.Cells.Item(97, 18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(97, 19).Value2 = 2000
$ws.Cells.Item(97, 20).Value2 = 2

# Ensure the D column (Fecha) keeps the date number format used throughout the column
$ws.Cells.Item(97, 4).NumberFormat = $ws.Cells.Item(98, 4).NumberFormat
